{"js": "// Remove the trailing generated-site footer block:\n//   - the blank paragraph right before \"Ver no Jupiter ...\"\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ...\"\n// while leaving the two pageBreakBefore marker paragraphs (and the blank\n// paragraph that follows the removed block) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst COPYRIGHT_PREFIX = \"\u00a9 2020 . Contact: luizeleno@usp.br.\";\nconst JUPITER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n\n// Locate the \"Ver no Jupiter ...\" paragraph, then remove it, the copyright\n// paragraph right after it, and the blank paragraph right before it.\nlet jupiterIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === JUPITER_TEXT) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex !== -1) {\n  const copyrightIndex = jupiterIndex + 1;\n  const blankIndex = jupiterIndex - 1;\n\n  if (\n    copyrightIndex < items.length &&\n    items[copyrightIndex].text.indexOf(COPYRIGHT_PREFIX) === 0\n  ) {\n    items[copyrightIndex].delete();\n  }\n\n  items[jupiterIndex].delete();\n\n  if (blankIndex >= 0 && items[blankIndex].text === \"\") {\n    items[blankIndex].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing generated-site footer block from the document:\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ...\"\n#   - the blank paragraph that sits right before \"Ver no Jupiter ...\"\n# The two pageBreakBefore marker (blank) paragraphs that bracket this block\n# are left untouched, as is the blank paragraph that follows the block.\n\n$d = $word.ActiveDocument\n$wdParagraph = 4\n\n# Locate + delete the copyright paragraph first (it comes after \"Ver no\n# Jupiter ...\" in the document, so removing it first does not shift the\n# position of the text we still need to find).\n$copyrightRange = $d.Content\n$foundCopyright = $copyrightRange.Find.Execute(\"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\")\nif ($foundCopyright) {\n    $copyrightRange.Expand($wdParagraph)\n    $copyrightRange.Delete()\n}\n\n# Locate the \"Ver no Jupiter ...\" paragraph.\n$jupiterRange = $d.Content\n$foundJupiter = $jupiterRange.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif ($foundJupiter) {\n    $jupiterRange.Expand($wdParagraph)\n\n    # Work out this paragraph's 1-based index so we can also reach the\n    # (empty) paragraph immediately before it.\n    $precedingRange = $d.Range(0, $jupiterRange.Start)\n    $jupiterIndex = $precedingRange.Paragraphs.Count + 1\n    $blankIndex = $jupiterIndex - 1\n\n    $blankText = \"\"\n    if ($blankIndex -ge 1) {\n        $blankText = $d.Paragraphs.Item($blankIndex).Range.Text.Trim()\n    }\n\n    $jupiterRange.Delete()\n\n    if ($blankIndex -ge 1 -and $blankText.Length -eq 0) {\n        $d.Paragraphs.Item($blankIndex).Range.Delete()\n    }\n}\n"}
